$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Hello, I'm Dr.Lee"
$ws.Range("B2").Value = "Hi, Dr. Lee."

$ws.Range("A3").Value = "What's your name?"
$ws.Range("B3").Value = "I'm Jin-ah."

$ws.Range("A4").Value = "Okay, Jin-ah. What brings you here?"
$ws.Range("B4").Value = "I've been feeling really down for a long time. It's not getting any better, and lately, I can't seem to sleep either. I just... I don't know what to do anymore."

$ws.Range("A5").Value = "Tell me about you more"
$ws.Range("B5").Value = "Well, I'm 24 years old, and I'm a graduate student. I live in Seoul with my family. I have a younger brother. My mom has had her own struggles with depression, and my dad has hypertension. I've been dealing with these feelings of sadness since early 2021, and it's just been getting worse. I started having trouble sleeping around mid-2022, and now it's almost impossible to get a good night's rest. I'm on antidepressants, but they don't seem to be helping much."
